$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.773.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.391.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.398.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0969'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.70%  '
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("E13").Value = '  -4.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.816.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.671.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.402.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '308.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.374'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.89%  '
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0721'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.54%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.814'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '133.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.564'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '249.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.92%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0906'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0486'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0210'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.88%  '
